$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 192.30027
$ws.Range("H2").Value = 576.90081
$ws.Range("I2").Value = 0.7642743413703218
$ws.Range("J2").Value = 0.7642743413703219
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 22.07422366666666
$ws.Range("N2").Value = 66.22267099999999
$ws.Range("O2").Value = 0.8730897844203874
$ws.Range("P2").Value = 0.8730897844203874
$ws.Range("Q2").Value = 4244.879171140389
$ws.Range("R2").Value = 38203.9125402635
$ws.Range("S2").Value = 0.6672801199450479
$ws.Range("T2").Value = 0.667280119945048
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 192.30027
$ws.Range("H3").Value = 576.90081
$ws.Range("I3").Value = 0.7642743413703218
$ws.Range("J3").Value = 0.7642743413703219
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.456833333333333
$ws.Range("N3").Value = 4.3705
$ws.Range("O3").Value = 0.05762133790721463
$ws.Range("P3").Value = 0.05762133790721465
$ws.Range("Q3").Value = 280.149443345
$ws.Range("R3").Value = 2521.344990105
$ws.Range("S3").Value = 0.04403851007791322
$ws.Range("T3").Value = 0.04403851007791323
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 192.30027
$ws.Range("H4").Value = 576.90081
$ws.Range("I4").Value = 0.7642743413703218
$ws.Range("J4").Value = 0.7642743413703219
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.751822333333333
$ws.Range("N4").Value = 5.255467
$ws.Range("O4").Value = 0.06928887767239802
$ws.Range("P4").Value = 0.06928887767239804
$ws.Range("Q4").Value = 336.87590769203
$ws.Range("R4").Value = 3031.88316922827
$ws.Range("S4").Value = 0.0529557113473608
$ws.Range("T4").Value = 0.05295571134736082
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 53.08542
$ws.Range("H5").Value = 159.25626
$ws.Range("I5").Value = 0.210981629962698
$ws.Range("J5").Value = 0.210981629962698
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 22.07422366666666
$ws.Range("N5").Value = 66.22267099999999
$ws.Range("O5").Value = 0.8730897844203874
$ws.Range("P5").Value = 0.8730897844203874
$ws.Range("Q5").Value = 1171.81943451894
$ws.Range("R5").Value = 10546.37491067046
$ws.Range("S5").Value = 0.1842059058207939
$ws.Range("T5").Value = 0.1842059058207939
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 53.08542
$ws.Range("H6").Value = 159.25626
$ws.Range("I6").Value = 0.210981629962698
$ws.Range("J6").Value = 0.210981629962698
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.456833333333333
$ws.Range("N6").Value = 4.3705
$ws.Range("O6").Value = 0.05762133790721463
$ws.Range("P6").Value = 0.05762133790721465
$ws.Range("Q6").Value = 77.33660936999999
$ws.Range("R6").Value = 696.0294843299999
$ws.Range("S6").Value = 0.01215704379229554
$ws.Range("T6").Value = 0.01215704379229554
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 53.08542
$ws.Range("H7").Value = 159.25626
$ws.Range("I7").Value = 0.210981629962698
$ws.Range("J7").Value = 0.210981629962698
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.751822333333333
$ws.Range("N7").Value = 5.255467
$ws.Range("O7").Value = 0.06928887767239802
$ws.Range("P7").Value = 0.06928887767239804
$ws.Range("Q7").Value = 92.99622433038
$ws.Range("R7").Value = 836.9660189734201
$ws.Range("S7").Value = 0.01461868034960853
$ws.Range("T7").Value = 0.01461868034960853
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 6.225884000000001
$ws.Range("H8").Value = 18.677652
$ws.Range("I8").Value = 0.02474402866698016
$ws.Range("J8").Value = 0.02474402866698016
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 22.07422366666666
$ws.Range("N8").Value = 66.22267099999999
$ws.Range("O8").Value = 0.8730897844203874
$ws.Range("P8").Value = 0.8730897844203874
$ws.Range("Q8").Value = 137.4315559387213
$ws.Range("R8").Value = 1236.884003448492
$ws.Range("S8").Value = 0.0216037586545456
$ws.Range("T8").Value = 0.0216037586545456
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 6.225884000000001
$ws.Range("H9").Value = 18.677652
$ws.Range("I9").Value = 0.02474402866698016
$ws.Range("J9").Value = 0.02474402866698016
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.456833333333333
$ws.Range("N9").Value = 4.3705
$ws.Range("O9").Value = 0.05762133790721463
$ws.Range("P9").Value = 0.05762133790721465
$ws.Range("Q9").Value = 9.070075340666667
$ws.Range("R9").Value = 81.630678066
$ws.Range("S9").Value = 0.00142578403700587
$ws.Range("T9").Value = 0.00142578403700587
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 6.225884000000001
$ws.Range("H10").Value = 18.677652
$ws.Range("I10").Value = 0.02474402866698016
$ws.Range("J10").Value = 0.02474402866698016
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.751822333333333
$ws.Range("N10").Value = 5.255467
$ws.Range("O10").Value = 0.06928887767239802
$ws.Range("P10").Value = 0.06928887767239804
$ws.Range("Q10").Value = 10.90664263594267
$ws.Range("R10").Value = 98.15978372348401
$ws.Range("S10").Value = 0.001714485975428699
$ws.Range("T10").Value = 0.001714485975428699
